# Auto-generated Excel COM-interop script.
# Adds new scrim-result rows to five worksheets, matching the
# commit "Actualizacion automatica de scrims_actualizado.xlsx (2025-07-27 04:04:40)".

$wb = $excel.ActiveWorkbook
# ===== Sheet: Belle's Rock =====
$ws = $wb.Worksheets.Item("Belle's Rock")

# Copy full-row formatting (A-F, H-N) from the last existing row (24) into the new rows
$ws.Range("A24:N24").Copy()
$ws.Range("A25:N26").PasteSpecial(-4122)

# Fix up the "G" (Equipo) column formatting per-row since it depends on the team value
$ws.Range("G24").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G26").PasteSpecial(-4122)

# Set the cell values for the new rows
$ws.Range("A25").Value = "BONNIE"
$ws.Range("B25").Value = "R-T"
$ws.Range("C25").Value = "TICK"
$ws.Range("D25").Value = "GUS"
$ws.Range("E25").Value = "SPROUT"
$ws.Range("F25").Value = "ANGELO"
$ws.Range("G25").Value = "Equipo 1"
$ws.Range("H25").Value = "LOUD|FireCrow"
$ws.Range("I25").Value = "LOUD|KaioDog"
$ws.Range("J25").Value = "LOUD|Edinho"
$ws.Range("K25").Value = "Doritos🐉"
$ws.Range("L25").Value = "CASA|Mohtep"
$ws.Range("M25").Value = "CASA|Pekka"
$ws.Range("N25").Value = "20250727T014855.000Z"

$ws.Range("A26").Value = "BONNIE"
$ws.Range("B26").Value = "R-T"
$ws.Range("C26").Value = "TICK"
$ws.Range("D26").Value = "GUS"
$ws.Range("E26").Value = "SPROUT"
$ws.Range("F26").Value = "ANGELO"
$ws.Range("G26").Value = "Equipo 1"
$ws.Range("H26").Value = "LOUD|FireCrow"
$ws.Range("I26").Value = "LOUD|KaioDog"
$ws.Range("J26").Value = "LOUD|Edinho"
$ws.Range("K26").Value = "Doritos🐉"
$ws.Range("L26").Value = "CASA|Mohtep"
$ws.Range("M26").Value = "CASA|Pekka"
$ws.Range("N26").Value = "20250727T014611.000Z"


# ===== Sheet: Dueling Beetles =====
$ws = $wb.Worksheets.Item("Dueling Beetles")

# Copy full-row formatting (A-F, H-N) from the last existing row (28) into the new rows
$ws.Range("A28:N28").Copy()
$ws.Range("A29:N31").PasteSpecial(-4122)

# Fix up the "G" (Equipo) column formatting per-row since it depends on the team value
$ws.Range("G27").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G28").Copy()
$ws.Range("G31").PasteSpecial(-4122)

# Set the cell values for the new rows
$ws.Range("A29").Value = "BERRY"
$ws.Range("B29").Value = "KAZE"
$ws.Range("C29").Value = "MEG"
$ws.Range("D29").Value = "BARLEY"
$ws.Range("E29").Value = "KIT"
$ws.Range("F29").Value = "CHARLIE"
$ws.Range("G29").Value = "Equipo 1"
$ws.Range("H29").Value = "LOUD|Edinho"
$ws.Range("I29").Value = "LOUD|FireCrow"
$ws.Range("J29").Value = "LOUD|KaioDog"
$ws.Range("K29").Value = "CASA|Pekka"
$ws.Range("L29").Value = "Doritos🐉"
$ws.Range("M29").Value = "CASA|Mohtep"
$ws.Range("N29").Value = "20250727T013929.000Z"

$ws.Range("A30").Value = "BERRY"
$ws.Range("B30").Value = "KAZE"
$ws.Range("C30").Value = "MEG"
$ws.Range("D30").Value = "BARLEY"
$ws.Range("E30").Value = "KIT"
$ws.Range("F30").Value = "CHARLIE"
$ws.Range("G30").Value = "Equipo 1"
$ws.Range("H30").Value = "LOUD|Edinho"
$ws.Range("I30").Value = "LOUD|FireCrow"
$ws.Range("J30").Value = "LOUD|KaioDog"
$ws.Range("K30").Value = "CASA|Pekka"
$ws.Range("L30").Value = "Doritos🐉"
$ws.Range("M30").Value = "CASA|Mohtep"
$ws.Range("N30").Value = "20250727T013741.000Z"

$ws.Range("A31").Value = "BERRY"
$ws.Range("B31").Value = "KAZE"
$ws.Range("C31").Value = "MEG"
$ws.Range("D31").Value = "BARLEY"
$ws.Range("E31").Value = "KIT"
$ws.Range("F31").Value = "CHARLIE"
$ws.Range("G31").Value = "Equipo 2"
$ws.Range("H31").Value = "LOUD|Edinho"
$ws.Range("I31").Value = "LOUD|FireCrow"
$ws.Range("J31").Value = "LOUD|KaioDog"
$ws.Range("K31").Value = "CASA|Pekka"
$ws.Range("L31").Value = "Doritos🐉"
$ws.Range("M31").Value = "CASA|Mohtep"
$ws.Range("N31").Value = "20250727T013540.000Z"


# ===== Sheet: Goldarm Gulch =====
$ws = $wb.Worksheets.Item("Goldarm Gulch")

# Copy full-row formatting (A-F, H-N) from the last existing row (20) into the new rows
$ws.Range("A20:N20").Copy()
$ws.Range("A21:N22").PasteSpecial(-4122)

# Fix up the "G" (Equipo) column formatting per-row since it depends on the team value
$ws.Range("G16").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("G21").PasteSpecial(-4122)

# Set the cell values for the new rows
$ws.Range("A21").Value = "LILY"
$ws.Range("B21").Value = "MANDY"
$ws.Range("C21").Value = "GENE"
$ws.Range("D21").Value = "BELLE"
$ws.Range("E21").Value = "BROCK"
$ws.Range("F21").Value = "DARRYL"
$ws.Range("G21").Value = "Equipo 2"
$ws.Range("H21").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I21").Value = "YT:BrabaoBS🎩"
$ws.Range("J21").Value = "Golden💘Mai"
$ws.Range("K21").Value = "FX|CaueBr"
$ws.Range("L21").Value = "FX|Wesley"
$ws.Range("M21").Value = "FX|REI DO FUT"
$ws.Range("N21").Value = "20250727T020141.000Z"

$ws.Range("A22").Value = "LILY"
$ws.Range("B22").Value = "MANDY"
$ws.Range("C22").Value = "GENE"
$ws.Range("D22").Value = "BELLE"
$ws.Range("E22").Value = "BROCK"
$ws.Range("F22").Value = "DARRYL"
$ws.Range("G22").Value = "Equipo 1"
$ws.Range("H22").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I22").Value = "YT:BrabaoBS🎩"
$ws.Range("J22").Value = "Golden💘Mai"
$ws.Range("K22").Value = "FX|CaueBr"
$ws.Range("L22").Value = "FX|Wesley"
$ws.Range("M22").Value = "FX|REI DO FUT"
$ws.Range("N22").Value = "20250727T015855.000Z"


# ===== Sheet: Double Swoosh =====
$ws = $wb.Worksheets.Item("Double Swoosh")

# Copy full-row formatting (A-F, H-N) from the last existing row (23) into the new rows
$ws.Range("A23:N23").Copy()
$ws.Range("A24:N26").PasteSpecial(-4122)

# Fix up the "G" (Equipo) column formatting per-row since it depends on the team value
$ws.Range("G23").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G20").Copy()
$ws.Range("G26").PasteSpecial(-4122)

# Set the cell values for the new rows
$ws.Range("A24").Value = "HANK"
$ws.Range("B24").Value = "CHESTER"
$ws.Range("C24").Value = "LUMI"
$ws.Range("D24").Value = "KAZE"
$ws.Range("E24").Value = "JANET"
$ws.Range("F24").Value = "BUSTER"
$ws.Range("G24").Value = "Equipo 1"
$ws.Range("H24").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I24").Value = "YT:BrabaoBS🎩"
$ws.Range("J24").Value = "Golden💘Mai"
$ws.Range("K24").Value = "FX|CaueBr"
$ws.Range("L24").Value = "FX|Wesley"
$ws.Range("M24").Value = "FX|REI DO FUT"
$ws.Range("N24").Value = "20250727T015224.000Z"

$ws.Range("A25").Value = "HANK"
$ws.Range("B25").Value = "CHESTER"
$ws.Range("C25").Value = "LUMI"
$ws.Range("D25").Value = "KAZE"
$ws.Range("E25").Value = "JANET"
$ws.Range("F25").Value = "BUSTER"
$ws.Range("G25").Value = "Equipo 1"
$ws.Range("H25").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I25").Value = "YT:BrabaoBS🎩"
$ws.Range("J25").Value = "Golden💘Mai"
$ws.Range("K25").Value = "FX|CaueBr"
$ws.Range("L25").Value = "FX|Wesley"
$ws.Range("M25").Value = "FX|REI DO FUT"
$ws.Range("N25").Value = "20250727T015006.000Z"

$ws.Range("A26").Value = "HANK"
$ws.Range("B26").Value = "CHESTER"
$ws.Range("C26").Value = "LUMI"
$ws.Range("D26").Value = "KAZE"
$ws.Range("E26").Value = "JANET"
$ws.Range("F26").Value = "BUSTER"
$ws.Range("G26").Value = "Equipo 2"
$ws.Range("H26").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I26").Value = "YT:BrabaoBS🎩"
$ws.Range("J26").Value = "Golden💘Mai"
$ws.Range("K26").Value = "FX|CaueBr"
$ws.Range("L26").Value = "FX|Wesley"
$ws.Range("M26").Value = "FX|REI DO FUT"
$ws.Range("N26").Value = "20250727T014710.000Z"


# ===== Sheet: Dry Season =====
$ws = $wb.Worksheets.Item("Dry Season")

# Copy full-row formatting (A-F, H-N) from the last existing row (65) into the new rows
$ws.Range("A65:N65").Copy()
$ws.Range("A66:N68").PasteSpecial(-4122)

# Fix up the "G" (Equipo) column formatting per-row since it depends on the team value
$ws.Range("G63").Copy()
$ws.Range("G66").PasteSpecial(-4122)
$ws.Range("G68").PasteSpecial(-4122)
$ws.Range("G65").Copy()
$ws.Range("G67").PasteSpecial(-4122)

# Set the cell values for the new rows
$ws.Range("A66").Value = "CARL"
$ws.Range("B66").Value = "BELLE"
$ws.Range("C66").Value = "BONNIE"
$ws.Range("D66").Value = "GUS"
$ws.Range("E66").Value = "BROCK"
$ws.Range("F66").Value = "R-T"
$ws.Range("G66").Value = "Equipo 1"
$ws.Range("H66").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I66").Value = "YT:BrabaoBS🎩"
$ws.Range("J66").Value = "Golden💘Mai"
$ws.Range("K66").Value = "FX|CaueBr"
$ws.Range("L66").Value = "FX|REI DO FUT"
$ws.Range("M66").Value = "FX|Wesley"
$ws.Range("N66").Value = "20250727T013955.000Z"

$ws.Range("A67").Value = "CARL"
$ws.Range("B67").Value = "BELLE"
$ws.Range("C67").Value = "BONNIE"
$ws.Range("D67").Value = "GUS"
$ws.Range("E67").Value = "BROCK"
$ws.Range("F67").Value = "R-T"
$ws.Range("G67").Value = "Equipo 2"
$ws.Range("H67").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I67").Value = "YT:BrabaoBS🎩"
$ws.Range("J67").Value = "Golden💘Mai"
$ws.Range("K67").Value = "FX|CaueBr"
$ws.Range("L67").Value = "FX|REI DO FUT"
$ws.Range("M67").Value = "FX|Wesley"
$ws.Range("N67").Value = "20250727T013739.000Z"

$ws.Range("A68").Value = "CARL"
$ws.Range("B68").Value = "BELLE"
$ws.Range("C68").Value = "BONNIE"
$ws.Range("D68").Value = "GUS"
$ws.Range("E68").Value = "BROCK"
$ws.Range("F68").Value = "R-T"
$ws.Range("G68").Value = "Equipo 1"
$ws.Range("H68").Value = "BKB|❄️IC€CRØW❄️"
$ws.Range("I68").Value = "YT:BrabaoBS🎩"
$ws.Range("J68").Value = "Golden💘Mai"
$ws.Range("K68").Value = "FX|CaueBr"
$ws.Range("L68").Value = "FX|REI DO FUT"
$ws.Range("M68").Value = "FX|Wesley"
$ws.Range("N68").Value = "20250727T013519.000Z"

